$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric, to avoid Excel auto-converting them to numbers
$textForceCells = @("D5", "D8", "D10", "D13", "D15", "D19", "D20", "D24", "D26", "D27", "D28", "D35", "D36", "D39", "D40", "D41", "D44", "D46", "D47", "D48", "D49")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range("D2").Value = "37.538.25"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.069.83"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "231.38"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "57.83"
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "2.376.00"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "14.75"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").Value = "0.763"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "2.081.17"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "37.481.33"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "6.16"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "69.89"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("D26").Value = "9.89"
$ws.Range("E26").Value = "  +4.69%  "
$ws.Range("D27").Value = "169.51"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").Value = "0.130"
$ws.Range("E28").Value = "  -5.26%  "
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("E30").Value = "  -4.87%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "2.54"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").Value = "1.83"
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "5.31"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").Value = "0.0226"
$ws.Range("E40").Value = "  +3.37%  "
$ws.Range("D41").Value = "98.20"
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.489.29"
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "2.90"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("D46").Value = "16.62"
$ws.Range("E46").Value = "  -3.66%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "1.04"
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "4.00"
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("D49").Value = "7.24"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").Value = "2.260.11"
$ws.Range("E51").Value = "  -0.77%  "

# Restore original (unstyled) appearance for the cells we temporarily forced to text format
foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}
